$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text columns D (price) to stay text even when they look numeric:
# set to Text format, write values, then ClearFormats to drop the temporary
# number-format style (keeps the cell style identical to its original state).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.456.13"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "3.501.76"
$ws.Range("E3").Value = "  -2.86%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "584.29"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "182.71"
$ws.Range("E6").Value = "  -3.96%  "
$ws.Range("D7").Value = "3.490.57"
$ws.Range("E7").Value = "  -3.11%  "
$ws.Range("D8").Value = "0.610"
$ws.Range("E8").Value = "  -3.48%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "0.196"
$ws.Range("E10").Value = "  +5.13%  "
$ws.Range("D11").Value = "0.641"
$ws.Range("E11").Value = "  -3.38%  "
$ws.Range("D12").Value = "53.82"
$ws.Range("E12").Value = "  -4.38%  "
$ws.Range("D13").Value = "0.0000303"
$ws.Range("E13").Value = "  -3.52%  "
$ws.Range("D14").Value = "9.40"
$ws.Range("E14").Value = "  -3.37%  "
$ws.Range("D15").Value = "4.057.94"
$ws.Range("E15").Value = "  -3.09%  "
$ws.Range("D16").Value = "19.21"
$ws.Range("E16").Value = "  -3.19%  "
$ws.Range("D17").Value = "69.400.94"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "3.476.47"
$ws.Range("E18").Value = "  -3.63%  "
$ws.Range("D19").Value = "12.27"
$ws.Range("E19").Value = "  -3.17%  "
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").Value = "532.15"
$ws.Range("E21").Value = "  +8.36%  "
$ws.Range("D22").Value = "1.01"
$ws.Range("E22").Value = "  -4.11%  "
$ws.Range("D23").Value = "18.39"
$ws.Range("E23").Value = "  -5.29%  "
$ws.Range("D24").Value = "4.54"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("D25").Value = "4.85"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "95.37"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "2.95"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "11.00"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").Value = "9.06"
$ws.Range("E29").Value = "  -4.01%  "
$ws.Range("D30").Value = "32.04"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").Value = "  -4.50%  "
$ws.Range("D32").Value = "12.37"
$ws.Range("D33").Value = "63.70"
$ws.Range("E33").Value = "  -3.44%  "
$ws.Range("D34").Value = "0.113"
$ws.Range("E34").Value = "  -4.59%  "
$ws.Range("D35").Value = "546.86"
$ws.Range("E35").Value = "  -6.07%  "
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "37.88"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "3.06"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "0.0₃0754"
$ws.Range("E40").Value = "  -7.97%  "
$ws.Range("E41").Value = "  -1.96%  "
$ws.Range("D42").Value = "3.347.97"
$ws.Range("E42").Value = "  +3.79%  "
$ws.Range("D43").Value = "3.36"
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "3.04"
$ws.Range("E44").Value = "  -6.95%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.49"
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("D46").Value = "2.95"
$ws.Range("E46").Value = "  -3.68%  "
$ws.Range("D47").Value = "0.0435"
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "8.99"
$ws.Range("E48").Value = "  -7.41%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "0.134"
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "137.06"
$ws.Range("E51").Value = "  +0.88%  "

$ws.Range("D2:D51").ClearFormats()
